$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (same layout as existing rows: A=Date, B=Method, C=ElapsedMs,
# D=wordCount, E=sentenceCount, F=posWordCount, G=negWordCount, H=posWordPercentage,
# I=negWordPercentage, J=positivePhraseCount, K=negativePhraseCount, L=posPhrasePercentage,
# M=negPhrasePercentage)
$rows = @(
    @(42600.786805555559, "Named", 5840, 2179, 111, 37, 9, 80, 19, 0, 0, 0, 0),
    @(42600.825462962966, "Named", 4954, 1423, 73, 18, 4, 81, 18, 0, 0, 0, 0)
)

$templateRow = $ws.Range("A2:M2")

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy formatting (including the date number format/style) from the existing data row
    $destRow = $ws.Range("A" + $r + ":M" + $r)
    $templateRow.Copy($destRow)

    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}

# Column A width grew slightly (auto fit due to new, slightly wider date text)
$ws.Columns.Item(1).ColumnWidth = 14
